# Update the mapping workbook:
#  - Bump the "Date" metadata value on the Metadata sheet
#  - Insert a new mapping row for FRCDAVaccination.doseQuantity <->
#    FRImmunizationDocument.doseQuantity on "Mapping Table 1" (right
#    after the "FRCDAVaccination.routeCode" row, before the
#    "consumable" row), shifting the rows below it down by one.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh the generation Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2026-01-22T09:24:45+00:00"

# --- Mapping Table 1: insert the new doseQuantity mapping row ---
$ws = $wb.Worksheets.Item("Mapping Table 1")

# Row 11 currently holds the "consumable" mapping; push it (and
# everything after it) down by inserting a fresh row above it.
$ws.Rows.Item(11).Insert()

# Re-apply the table's row style (border/fill/alignment) to the newly
# inserted row by copying it from the row right below (same formatting
# as every other data row) before filling in the new values.
$ws.Range("A12:E12").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "FRCDAVaccination.doseQuantity"
$ws.Cells.Item(11, 3).Value = "equivalent"
$ws.Cells.Item(11, 4).Value = "FRImmunizationDocument.doseQuantity"
